# The "Förändrad" (changed) date in column C was bumped by one day
# (2023-10-04 -> 2023-10-05, serial 45203 -> 45204) for every data row
# (rows 2 through 307) on the single worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C307").Value = 45204
